# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Map of sheet name -> hashtable of row -> new F-column value
$updates = @{
    "展览"   = @{ 2 = 6916; 7 = 549; 8 = 124; 13 = 192; 18 = 3528; 22 = 2152; 23 = 213 }
    "全部类型" = @{ 2 = 6916; 8 = 549; 9 = 124; 14 = 192; 19 = 3528; 23 = 2152; 24 = 213 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
